$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 7).Value = 27.524977
$ws.Cells.Item(2, 8).Value = 55.049954
$ws.Cells.Item(2, 9).Value = 0.02927501708753065
$ws.Cells.Item(2, 10).Value = 0.02011213558514335
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 13).Value = 19.5719925
$ws.Cells.Item(2, 14).Value = 39.143985
$ws.Cells.Item(2, 15).Value = 0.03094210933382397
$ws.Cells.Item(2, 16).Value = 0.02187190777676379
$ws.Cells.Item(2, 17).Value = 538.7186434066725
$ws.Cells.Item(2, 18).Value = 2154.87457362669
$ws.Cells.Item(2, 19).Value = 0.0009058307794719384
$ws.Cells.Item(2, 20).Value = 0.0004398907747120246

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 7).Value = 27.524977
$ws.Cells.Item(3, 8).Value = 55.049954
$ws.Cells.Item(3, 9).Value = 0.02927501708753065
$ws.Cells.Item(3, 10).Value = 0.02011213558514335
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 115.495743
$ws.Cells.Item(3, 14).Value = 346.487229
$ws.Cells.Item(3, 15).Value = 0.1825916246134488
$ws.Cells.Item(3, 16).Value = 0.1936015640337701
$ws.Cells.Item(3, 17).Value = 3179.017669672911
$ws.Cells.Item(3, 18).Value = 19074.10601803747
$ws.Cells.Item(3, 19).Value = 0.005345372930598695
$ws.Cells.Item(3, 20).Value = 0.003893740905342997

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 7).Value = 27.524977
$ws.Cells.Item(4, 8).Value = 55.049954
$ws.Cells.Item(4, 9).Value = 0.02927501708753065
$ws.Cells.Item(4, 10).Value = 0.02011213558514335
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 239.8982746666667
$ws.Cells.Item(4, 14).Value = 719.694824
$ws.Cells.Item(4, 15).Value = 0.3792643310961689
$ws.Cells.Item(4, 16).Value = 0.4021332732970914
$ws.Cells.Item(4, 17).Value = 6603.194492539683
$ws.Cells.Item(4, 18).Value = 39619.1669552381
$ws.Cells.Item(4, 19).Value = 0.01110296977353123
$ws.Cells.Item(4, 20).Value = 0.00808775891584861

$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 7).Value = 27.524977
$ws.Cells.Item(5, 8).Value = 55.049954
$ws.Cells.Item(5, 9).Value = 0.02927501708753065
$ws.Cells.Item(5, 10).Value = 0.02011213558514335
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 151.102183
$ws.Cells.Item(5, 14).Value = 453.306549
$ws.Cells.Item(5, 15).Value = 0.2388832034840335
$ws.Cells.Item(5, 16).Value = 0.2532874216646837
$ws.Cells.Item(5, 17).Value = 4159.084111724791
$ws.Cells.Item(5, 18).Value = 24954.50467034874
$ws.Cells.Item(5, 19).Value = 0.006993309863919142
$ws.Cells.Item(5, 20).Value = 0.005094150966531495

$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 7).Value = 27.524977
$ws.Cells.Item(6, 8).Value = 55.049954
$ws.Cells.Item(6, 9).Value = 0.02927501708753065
$ws.Cells.Item(6, 10).Value = 0.02011213558514335
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 18.12446233333333
$ws.Cells.Item(6, 14).Value = 54.373387
$ws.Cells.Item(6, 15).Value = 0.02865365369084289
$ws.Cells.Item(6, 16).Value = 0.03038141635232813
$ws.Cells.Item(6, 17).Value = 498.8754088623663
$ws.Cells.Item(6, 18).Value = 2993.252453174198
$ws.Cells.Item(6, 19).Value = 0.0008388362014196113
$ws.Cells.Item(6, 20).Value = 0.0006110351649467147

$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 7).Value = 27.524977
$ws.Cells.Item(7, 8).Value = 55.049954
$ws.Cells.Item(7, 9).Value = 0.02927501708753065
$ws.Cells.Item(7, 10).Value = 0.02011213558514335
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 13).Value = 88.3431645
$ws.Cells.Item(7, 14).Value = 176.686329
$ws.Cells.Item(7, 15).Value = 0.139665077781682
$ws.Cells.Item(7, 16).Value = 0.09872441687536272
$ws.Cells.Item(7, 17).Value = 2431.643570969717
$ws.Cells.Item(7, 18).Value = 9726.574283878866
$ws.Cells.Item(7, 19).Value = 0.004088697538590038
$ws.Cells.Item(7, 20).Value = 0.00198555885776151

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 830.1578366666666
$ws.Cells.Item(8, 8).Value = 2490.47351
$ws.Cells.Item(8, 9).Value = 0.8829393337463696
$ws.Cells.Item(8, 10).Value = 0.9098779792682091
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 13).Value = 19.5719925
$ws.Cells.Item(8, 14).Value = 39.143985
$ws.Cells.Item(8, 15).Value = 0.03094210933382397
$ws.Cells.Item(8, 16).Value = 0.02187190777676379
$ws.Cells.Item(8, 17).Value = 16247.84295305622
$ws.Cells.Item(8, 18).Value = 97487.05771833734
$ws.Cells.Item(8, 19).Value = 0.02732000539991386
$ws.Cells.Item(8, 20).Value = 0.01990076725066246

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 830.1578366666666
$ws.Cells.Item(9, 8).Value = 2490.47351
$ws.Cells.Item(9, 9).Value = 0.8829393337463696
$ws.Cells.Item(9, 10).Value = 0.9098779792682091
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 115.495743
$ws.Cells.Item(9, 14).Value = 346.487229
$ws.Cells.Item(9, 15).Value = 0.1825916246134488
$ws.Cells.Item(9, 16).Value = 0.1936015640337701
$ws.Cells.Item(9, 17).Value = 95879.69615308932
$ws.Cells.Item(9, 18).Value = 862917.2653778037
$ws.Cells.Item(9, 19).Value = 0.1612173273838657
$ws.Cells.Item(9, 20).Value = 0.1761537998662115

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 830.1578366666666
$ws.Cells.Item(10, 8).Value = 2490.47351
$ws.Cells.Item(10, 9).Value = 0.8829393337463696
$ws.Cells.Item(10, 10).Value = 0.9098779792682091
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 239.8982746666667
$ws.Cells.Item(10, 14).Value = 719.694824
$ws.Cells.Item(10, 15).Value = 0.3792643310961689
$ws.Cells.Item(10, 16).Value = 0.4021332732970914
$ws.Cells.Item(10, 17).Value = 199153.4327173458
$ws.Cells.Item(10, 18).Value = 1792380.894456112
$ws.Cells.Item(10, 19).Value = 0.3348673958118139
$ws.Cells.Item(10, 20).Value = 0.365892210104068

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 830.1578366666666
$ws.Cells.Item(11, 8).Value = 2490.47351
$ws.Cells.Item(11, 9).Value = 0.8829393337463696
$ws.Cells.Item(11, 10).Value = 0.9098779792682091
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 151.102183
$ws.Cells.Item(11, 14).Value = 453.306549
$ws.Cells.Item(11, 15).Value = 0.2388832034840335
$ws.Cells.Item(11, 16).Value = 0.2532874216646837
$ws.Cells.Item(11, 17).Value = 125438.6613548908
$ws.Cells.Item(11, 18).Value = 1128947.952194017
$ws.Cells.Item(11, 19).Value = 0.210919376527391
$ws.Cells.Item(11, 20).Value = 0.2304606473983172

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 830.1578366666666
$ws.Cells.Item(12, 8).Value = 2490.47351
$ws.Cells.Item(12, 9).Value = 0.8829393337463696
$ws.Cells.Item(12, 10).Value = 0.9098779792682091
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 18.12446233333333
$ws.Cells.Item(12, 14).Value = 54.373387
$ws.Cells.Item(12, 15).Value = 0.02865365369084289
$ws.Cells.Item(12, 16).Value = 0.03038141635232813
$ws.Cells.Item(12, 17).Value = 15046.16444138649
$ws.Cells.Item(12, 18).Value = 135415.4799724784
$ws.Cells.Item(12, 19).Value = 0.02529943789919202
$ws.Cells.Item(12, 20).Value = 0.02764338171796244

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 830.1578366666666
$ws.Cells.Item(13, 8).Value = 2490.47351
$ws.Cells.Item(13, 9).Value = 0.8829393337463696
$ws.Cells.Item(13, 10).Value = 0.9098779792682091
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 13).Value = 88.3431645
$ws.Cells.Item(13, 14).Value = 176.686329
$ws.Cells.Item(13, 15).Value = 0.139665077781682
$ws.Cells.Item(13, 16).Value = 0.09872441687536272
$ws.Cells.Item(13, 17).Value = 73338.77032560746
$ws.Cells.Item(13, 18).Value = 440032.6219536447
$ws.Cells.Item(13, 19).Value = 0.1233157907241932
$ws.Cells.Item(13, 20).Value = 0.08982717293098731

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 1.012325
$ws.Cells.Item(14, 8).Value = 3.036975
$ws.Cells.Item(14, 9).Value = 0.001076688698890992
$ws.Cells.Item(14, 10).Value = 0.001109538674068479
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 13).Value = 19.5719925
$ws.Cells.Item(14, 14).Value = 39.143985
$ws.Cells.Item(14, 15).Value = 0.03094210933382397
$ws.Cells.Item(14, 16).Value = 0.02187190777676379
$ws.Cells.Item(14, 17).Value = 19.8132173075625
$ws.Cells.Item(14, 18).Value = 118.879303845375
$ws.Cells.Item(14, 19).Value = 0.00003331501943957774
$ws.Cells.Item(14, 20).Value = 0.00002426772755397854

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 1.012325
$ws.Cells.Item(15, 8).Value = 3.036975
$ws.Cells.Item(15, 9).Value = 0.001076688698890992
$ws.Cells.Item(15, 10).Value = 0.001109538674068479
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 115.495743
$ws.Cells.Item(15, 14).Value = 346.487229
$ws.Cells.Item(15, 15).Value = 0.1825916246134488
$ws.Cells.Item(15, 16).Value = 0.1936015640337701
$ws.Cells.Item(15, 17).Value = 116.919228032475
$ws.Cells.Item(15, 18).Value = 1052.273052292275
$ws.Cells.Item(15, 19).Value = 0.0001965943387334465
$ws.Cells.Item(15, 20).Value = 0.000214808422655613

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 1.012325
$ws.Cells.Item(16, 8).Value = 3.036975
$ws.Cells.Item(16, 9).Value = 0.001076688698890992
$ws.Cells.Item(16, 10).Value = 0.001109538674068479
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 239.8982746666667
$ws.Cells.Item(16, 14).Value = 719.694824
$ws.Cells.Item(16, 15).Value = 0.3792643310961689
$ws.Cells.Item(16, 16).Value = 0.4021332732970914
$ws.Cells.Item(16, 17).Value = 242.8550209019334
$ws.Cells.Item(16, 18).Value = 2185.6951881174
$ws.Cells.Item(16, 19).Value = 0.0004083496191836963
$ws.Cells.Item(16, 20).Value = 0.000446182418852872

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 1.012325
$ws.Cells.Item(17, 8).Value = 3.036975
$ws.Cells.Item(17, 9).Value = 0.001076688698890992
$ws.Cells.Item(17, 10).Value = 0.001109538674068479
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 151.102183
$ws.Cells.Item(17, 14).Value = 453.306549
$ws.Cells.Item(17, 15).Value = 0.2388832034840335
$ws.Cells.Item(17, 16).Value = 0.2532874216646837
$ws.Cells.Item(17, 17).Value = 152.964517405475
$ws.Cells.Item(17, 18).Value = 1376.680656649275
$ws.Cells.Item(17, 19).Value = 0.000257202845546136
$ws.Cells.Item(17, 20).Value = 0.0002810321899920568

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 1.012325
$ws.Cells.Item(18, 8).Value = 3.036975
$ws.Cells.Item(18, 9).Value = 0.001076688698890992
$ws.Cells.Item(18, 10).Value = 0.001109538674068479
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 18.12446233333333
$ws.Cells.Item(18, 14).Value = 54.373387
$ws.Cells.Item(18, 15).Value = 0.02865365369084289
$ws.Cells.Item(18, 16).Value = 0.03038141635232813
$ws.Cells.Item(18, 17).Value = 18.34784633159167
$ws.Cells.Item(18, 18).Value = 165.130616984325
$ws.Cells.Item(18, 19).Value = 0.00003085106511086669
$ws.Cells.Item(18, 20).Value = 0.00003370935641588455

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 1.012325
$ws.Cells.Item(19, 8).Value = 3.036975
$ws.Cells.Item(19, 9).Value = 0.001076688698890992
$ws.Cells.Item(19, 10).Value = 0.001109538674068479
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 13).Value = 88.3431645
$ws.Cells.Item(19, 14).Value = 176.686329
$ws.Cells.Item(19, 15).Value = 0.139665077781682
$ws.Cells.Item(19, 16).Value = 0.09872441687536272
$ws.Cells.Item(19, 17).Value = 89.43199400246252
$ws.Cells.Item(19, 18).Value = 536.591964014775
$ws.Cells.Item(19, 19).Value = 0.0001503758108772683
$ws.Cells.Item(19, 20).Value = 0.0001095385585980737

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 1.320133
$ws.Cells.Item(20, 8).Value = 3.960399
$ws.Cells.Item(20, 9).Value = 0.001404067154454411
$ws.Cells.Item(20, 10).Value = 0.001446905508027602
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 13).Value = 19.5719925
$ws.Cells.Item(20, 14).Value = 39.143985
$ws.Cells.Item(20, 15).Value = 0.03094210933382397
$ws.Cells.Item(20, 16).Value = 0.02187190777676379
$ws.Cells.Item(20, 17).Value = 25.8376331750025
$ws.Cells.Item(20, 18).Value = 155.025799050015
$ws.Cells.Item(20, 19).Value = 0.00004344479940515949
$ws.Cells.Item(20, 20).Value = 0.00003164658383327128

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 1.320133
$ws.Cells.Item(21, 8).Value = 3.960399
$ws.Cells.Item(21, 9).Value = 0.001404067154454411
$ws.Cells.Item(21, 10).Value = 0.001446905508027602
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 115.495743
$ws.Cells.Item(21, 14).Value = 346.487229
$ws.Cells.Item(21, 15).Value = 0.1825916246134488
$ws.Cells.Item(21, 16).Value = 0.1936015640337701
$ws.Cells.Item(21, 17).Value = 152.469741693819
$ws.Cells.Item(21, 18).Value = 1372.227675244371
$ws.Cells.Item(21, 19).Value = 0.000256370902798213
$ws.Cells.Item(21, 20).Value = 0.0002801231693632206

$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 1.320133
$ws.Cells.Item(22, 8).Value = 3.960399
$ws.Cells.Item(22, 9).Value = 0.001404067154454411
$ws.Cells.Item(22, 10).Value = 0.001446905508027602
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 239.8982746666667
$ws.Cells.Item(22, 14).Value = 719.694824
$ws.Cells.Item(22, 15).Value = 0.3792643310961689
$ws.Cells.Item(22, 16).Value = 0.4021332732970914
$ws.Cells.Item(22, 17).Value = 316.6976290305307
$ws.Cells.Item(22, 18).Value = 2850.278661274776
$ws.Cells.Item(22, 19).Value = 0.0005325125901482534
$ws.Cells.Item(22, 20).Value = 0.0005818488480947308

$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 1.320133
$ws.Cells.Item(23, 8).Value = 3.960399
$ws.Cells.Item(23, 9).Value = 0.001404067154454411
$ws.Cells.Item(23, 10).Value = 0.001446905508027602
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 151.102183
$ws.Cells.Item(23, 14).Value = 453.306549
$ws.Cells.Item(23, 15).Value = 0.2388832034840335
$ws.Cells.Item(23, 16).Value = 0.2532874216646837
$ws.Cells.Item(23, 17).Value = 199.474978150339
$ws.Cells.Item(23, 18).Value = 1795.274803353051
$ws.Cells.Item(23, 19).Value = 0.0003354080597627809
$ws.Cells.Item(23, 20).Value = 0.0003664829655207407

$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 1.320133
$ws.Cells.Item(24, 8).Value = 3.960399
$ws.Cells.Item(24, 9).Value = 0.001404067154454411
$ws.Cells.Item(24, 10).Value = 0.001446905508027602
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 18.12446233333333
$ws.Cells.Item(24, 14).Value = 54.373387
$ws.Cells.Item(24, 15).Value = 0.02865365369084289
$ws.Cells.Item(24, 16).Value = 0.03038141635232813
$ws.Cells.Item(24, 17).Value = 23.92670083349033
$ws.Cells.Item(24, 18).Value = 215.340307501413
$ws.Cells.Item(24, 19).Value = 0.0000402316540024239
$ws.Cells.Item(24, 20).Value = 0.00004395903866186344

$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 1.320133
$ws.Cells.Item(25, 8).Value = 3.960399
$ws.Cells.Item(25, 9).Value = 0.001404067154454411
$ws.Cells.Item(25, 10).Value = 0.001446905508027602
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 13).Value = 88.3431645
$ws.Cells.Item(25, 14).Value = 176.686329
$ws.Cells.Item(25, 15).Value = 0.139665077781682
$ws.Cells.Item(25, 16).Value = 0.09872441687536272
$ws.Cells.Item(25, 17).Value = 116.6247267808785
$ws.Cells.Item(25, 18).Value = 699.748360685271
$ws.Cells.Item(25, 19).Value = 0.0001960991483375802
$ws.Cells.Item(25, 20).Value = 0.0001428449025537755

$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 24.21940733333333
$ws.Cells.Item(26, 8).Value = 72.658222
$ws.Cells.Item(26, 9).Value = 0.02575927905528127
$ws.Cells.Item(26, 10).Value = 0.02654519951532467
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 19.5719925
$ws.Cells.Item(26, 14).Value = 39.143985
$ws.Cells.Item(26, 15).Value = 0.03094210933382397
$ws.Cells.Item(26, 16).Value = 0.02187190777676379
$ws.Cells.Item(26, 17).Value = 474.022058682445
$ws.Cells.Item(26, 18).Value = 2844.13235209467
$ws.Cells.Item(26, 19).Value = 0.0007970464288889948
$ws.Cells.Item(26, 20).Value = 0.0005805941557149761

$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 7).Value = 24.21940733333333
$ws.Cells.Item(27, 8).Value = 72.658222
$ws.Cells.Item(27, 9).Value = 0.02575927905528127
$ws.Cells.Item(27, 10).Value = 0.02654519951532467
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 13).Value = 115.495743
$ws.Cells.Item(27, 14).Value = 346.487229
$ws.Cells.Item(27, 15).Value = 0.1825916246134488
$ws.Cells.Item(27, 16).Value = 0.1936015640337701
$ws.Cells.Item(27, 17).Value = 2797.238444982982
$ws.Cells.Item(27, 18).Value = 25175.14600484684
$ws.Cells.Item(27, 19).Value = 0.00470342861157499
$ws.Cells.Item(27, 20).Value = 0.005139192143755333

$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 7).Value = 24.21940733333333
$ws.Cells.Item(28, 8).Value = 72.658222
$ws.Cells.Item(28, 9).Value = 0.02575927905528127
$ws.Cells.Item(28, 10).Value = 0.02654519951532467
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 13).Value = 239.8982746666667
$ws.Cells.Item(28, 14).Value = 719.694824
$ws.Cells.Item(28, 15).Value = 0.3792643310961689
$ws.Cells.Item(28, 16).Value = 0.4021332732970914
$ws.Cells.Item(28, 17).Value = 5810.194032715881
$ws.Cells.Item(28, 18).Value = 52291.74629444293
$ws.Cells.Item(28, 19).Value = 0.009769575740420802
$ws.Cells.Item(28, 20).Value = 0.01067470797142188

$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 7).Value = 24.21940733333333
$ws.Cells.Item(29, 8).Value = 72.658222
$ws.Cells.Item(29, 9).Value = 0.02575927905528127
$ws.Cells.Item(29, 10).Value = 0.02654519951532467
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 13).Value = 151.102183
$ws.Cells.Item(29, 14).Value = 453.306549
$ws.Cells.Item(29, 15).Value = 0.2388832034840335
$ws.Cells.Item(29, 16).Value = 0.2532874216646837
$ws.Cells.Item(29, 17).Value = 3659.605319032875
$ws.Cells.Item(29, 18).Value = 32936.44787129587
$ws.Cells.Item(29, 19).Value = 0.006153459100164757
$ws.Cells.Item(29, 20).Value = 0.006723565142811199

$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 7).Value = 24.21940733333333
$ws.Cells.Item(30, 8).Value = 72.658222
$ws.Cells.Item(30, 9).Value = 0.02575927905528127
$ws.Cells.Item(30, 10).Value = 0.02654519951532467
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 13).Value = 18.12446233333333
$ws.Cells.Item(30, 14).Value = 54.373387
$ws.Cells.Item(30, 15).Value = 0.02865365369084289
$ws.Cells.Item(30, 16).Value = 0.03038141635232813
$ws.Cells.Item(30, 17).Value = 438.9637359486571
$ws.Cells.Item(30, 18).Value = 3950.673623537914
$ws.Cells.Item(30, 19).Value = 0.000738097461375812
$ws.Cells.Item(30, 20).Value = 0.0008064807586306978

$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 7).Value = 24.21940733333333
$ws.Cells.Item(31, 8).Value = 72.658222
$ws.Cells.Item(31, 9).Value = 0.02575927905528127
$ws.Cells.Item(31, 10).Value = 0.02654519951532467
$ws.Cells.Item(31, 11).Value = 2
$ws.Cells.Item(31, 13).Value = 88.3431645
$ws.Cells.Item(31, 14).Value = 176.686329
$ws.Cells.Item(31, 15).Value = 0.139665077781682
$ws.Cells.Item(31, 16).Value = 0.09872441687536272
$ws.Cells.Item(31, 17).Value = 2139.619086141173
$ws.Cells.Item(31, 18).Value = 12837.71451684704
$ws.Cells.Item(31, 19).Value = 0.003597671712855909
$ws.Cells.Item(31, 20).Value = 0.002620659342990589

$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 7).Value = 55.986019
$ws.Cells.Item(32, 8).Value = 111.972038
$ws.Cells.Item(32, 9).Value = 0.05954561425747298
$ws.Cells.Item(32, 10).Value = 0.04090824144922671
$ws.Cells.Item(32, 11).Value = 2
$ws.Cells.Item(32, 13).Value = 19.5719925
$ws.Cells.Item(32, 14).Value = 39.143985
$ws.Cells.Item(32, 15).Value = 0.03094210933382397
$ws.Cells.Item(32, 16).Value = 0.02187190777676379
$ws.Cells.Item(32, 17).Value = 1095.757943972857
$ws.Cells.Item(32, 18).Value = 4383.03177589143
$ws.Cells.Item(32, 19).Value = 0.001842466906704436
$ws.Cells.Item(32, 20).Value = 0.0008947412842870724

$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 7).Value = 55.986019
$ws.Cells.Item(33, 8).Value = 111.972038
$ws.Cells.Item(33, 9).Value = 0.05954561425747298
$ws.Cells.Item(33, 10).Value = 0.04090824144922671
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 13).Value = 115.495743
$ws.Cells.Item(33, 14).Value = 346.487229
$ws.Cells.Item(33, 15).Value = 0.1825916246134488
$ws.Cells.Item(33, 16).Value = 0.1936015640337701
$ws.Cells.Item(33, 17).Value = 6466.146862017117
$ws.Cells.Item(33, 18).Value = 38796.8811721027
$ws.Cells.Item(33, 19).Value = 0.01087253044587773
$ws.Cells.Item(33, 20).Value = 0.007919899526441393

$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(34, 7).Value = 55.986019
$ws.Cells.Item(34, 8).Value = 111.972038
$ws.Cells.Item(34, 9).Value = 0.05954561425747298
$ws.Cells.Item(34, 10).Value = 0.04090824144922671
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 13).Value = 239.8982746666667
$ws.Cells.Item(34, 14).Value = 719.694824
$ws.Cells.Item(34, 15).Value = 0.3792643310961689
$ws.Cells.Item(34, 16).Value = 0.4021332732970914
$ws.Cells.Item(34, 17).Value = 13430.94936355522
$ws.Cells.Item(34, 18).Value = 80585.69618133131
$ws.Cells.Item(34, 19).Value = 0.02258352756107098
$ws.Cells.Item(34, 20).Value = 0.01645056503880529

$ws.Cells.Item(35, 5).Value = 2
$ws.Cells.Item(35, 7).Value = 55.986019
$ws.Cells.Item(35, 8).Value = 111.972038
$ws.Cells.Item(35, 9).Value = 0.05954561425747298
$ws.Cells.Item(35, 10).Value = 0.04090824144922671
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 13).Value = 151.102183
$ws.Cells.Item(35, 14).Value = 453.306549
$ws.Cells.Item(35, 15).Value = 0.2388832034840335
$ws.Cells.Item(35, 16).Value = 0.2532874216646837
$ws.Cells.Item(35, 17).Value = 8459.609688379476
$ws.Cells.Item(35, 18).Value = 50757.65813027685
$ws.Cells.Item(35, 19).Value = 0.01422444708724968
$ws.Cells.Item(35, 20).Value = 0.01036154300151098

$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 7).Value = 55.986019
$ws.Cells.Item(36, 8).Value = 111.972038
$ws.Cells.Item(36, 9).Value = 0.05954561425747298
$ws.Cells.Item(36, 10).Value = 0.04090824144922671
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 13).Value = 18.12446233333333
$ws.Cells.Item(36, 14).Value = 54.373387
$ws.Cells.Item(36, 15).Value = 0.02865365369084289
$ws.Cells.Item(36, 16).Value = 0.03038141635232813
$ws.Cells.Item(36, 17).Value = 1014.716492558784
$ws.Cells.Item(36, 18).Value = 6088.298955352706
$ws.Cells.Item(36, 19).Value = 0.001706199409742147
$ws.Cells.Item(36, 20).Value = 0.001242850315710524

$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 7).Value = 55.986019
$ws.Cells.Item(37, 8).Value = 111.972038
$ws.Cells.Item(37, 9).Value = 0.05954561425747298
$ws.Cells.Item(37, 10).Value = 0.04090824144922671
$ws.Cells.Item(37, 11).Value = 2
$ws.Cells.Item(37, 13).Value = 88.3431645
$ws.Cells.Item(37, 14).Value = 176.686329
$ws.Cells.Item(37, 15).Value = 0.139665077781682
$ws.Cells.Item(37, 16).Value = 0.09872441687536272
$ws.Cells.Item(37, 17).Value = 4945.982086217125
$ws.Cells.Item(37, 18).Value = 19783.9283448685
$ws.Cells.Item(37, 19).Value = 0.008316442846827993
$ws.Cells.Item(37, 20).Value = 0.00403864228247145
